$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update prepends two new daily observations (rows 333-334) to the
# Brócoli series, pushing the existing rows 333-429 down to 335-431.
$ws.Rows("333:334").Insert()

# --- New row 333 --------------------------------------------------------
$ws.Range("A333").Value = 1
$ws.Range("B333").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C333").Value = "Arica y Parinacota"
$ws.Range("D333").Value = 44809
$ws.Range("E333").Value = 15
$ws.Range("F333").Value = 100112023
$ws.Range("G333").Value = "Brócoli"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Segunda"
$ws.Range("J333").Value = 900
$ws.Range("K333").Value = 700
$ws.Range("L333").Value = 800
$ws.Range("M333").Value = 750
$ws.Range("N333").Value = '$/unidad'
$ws.Range("O333").Value = "Región de Arica y Parinacota"
$ws.Range("P333").Value = 750
$ws.Range("Q333").Value = 1
$ws.Range("R333").Value = "Hortaliza"

# --- New row 334 --------------------------------------------------------
$ws.Range("A334").Value = 1
$ws.Range("B334").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C334").Value = "Arica y Parinacota"
$ws.Range("D334").Value = 44809
$ws.Range("E334").Value = 15
$ws.Range("F334").Value = 100112023
$ws.Range("G334").Value = "Brócoli"
$ws.Range("H334").Value = "Sin especificar"
$ws.Range("I334").Value = "Tercera"
$ws.Range("J334").Value = 800
$ws.Range("K334").Value = 500
$ws.Range("L334").Value = 600
$ws.Range("M334").Value = 550
$ws.Range("N334").Value = '$/unidad'
$ws.Range("O334").Value = "Región de Arica y Parinacota"
$ws.Range("P334").Value = 550
$ws.Range("Q334").Value = 1
$ws.Range("R334").Value = "Hortaliza"
